# Leave Card update — 12/27/2023 4:01 PM
# Adds new leave entries (SP, SL, VL) across several pay periods, extends the
# bi-weekly period start dates down through Dec 2024, and grows Table1 by
# three rows so later entries keep landing inside the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 119: SP(1-0-0) entry gets its EARNED value -------------------------
$ws.Range("C119").Value = 1.25

# --- Row 122: new SP(2-0-0) entry -------------------------------------------
$ws.Range("A122").Value = 45200
$ws.Range("B122").Value = "SP(2-0-0)"
$ws.Range("C122").Value = 1.25
$ws.Range("K122").Value = "10/5,6/2023"

# --- Rows 123-125: SL / VL entries ------------------------------------------
$ws.Range("B123").Value = "SL(1-0-0)"
$ws.Range("H123").Value = 1

$ws.Range("B124").Value = "SL(1-0-0)"
$ws.Range("H124").Value = 1

$ws.Range("B125").Value = "VL(1-0-0)"
$ws.Range("D125").Value = 1

# K123:K125 need the date number format used elsewhere in the REMARKS column
# (style 49) rather than the plain style already on those cells, so copy it
# over from an existing dated remark before writing the new values.
$ws.Range("K118").Copy()
$ws.Range("K123:K125").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("K123").Value = 45210
$ws.Range("K124").Value = 45223
$ws.Range("K125").Value = 45225

# --- Rows 126-140: bi-weekly period start dates continuing into 2024 -------
$ws.Range("A126").Value = 45231
$ws.Range("A127").Value = 45261

# Row 128 is a year-break row like row 103 ("2023"): copy that row's format
# (style 48, which quote-prefixes the value as text) then write "2024".
$ws.Range("A103").Copy()
$ws.Range("A128").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A128").Value = "'2024"

$ws.Range("A129").Value = 45292
$ws.Range("A130").Value = 45323
$ws.Range("A131").Value = 45352
$ws.Range("A132").Value = 45383
$ws.Range("A133").Value = 45413
$ws.Range("A134").Value = 45444
$ws.Range("A135").Value = 45474
$ws.Range("A136").Value = 45505
$ws.Range("A137").Value = 45536
$ws.Range("A138").Value = 45566
$ws.Range("A139").Value = 45597
$ws.Range("A140").Value = 45627

# --- Grow Table1 by three rows (149 -> 152) ---------------------------------
$lo = $ws.ListObjects.Item("Table1")

# Capture the current (still "last row") bottom-border formatting of row 149
# onto the new bottom row 152 before we touch row 149 itself.
$ws.Range("A149:K149").Copy()
$ws.Range("A152:K152").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Rows 149-151 become ordinary interior rows now, matching row 148's look.
$ws.Range("A148:K148").Copy()
$ws.Range("A149:K151").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Formats-only paste doesn't bring the calculated column's formula along, so
# restore it on every row the copies touched.
$earnedFormula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),`"`",Table1[[#This Row],[EARNED]])"
$ws.Range("G149").Formula = $earnedFormula
$ws.Range("G150").Formula = $earnedFormula
$ws.Range("G151").Formula = $earnedFormula
$ws.Range("G152").Formula = $earnedFormula

# --- Match the author's last selection before saving ------------------------
$ws.Range("K125").Select()
